$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, shifting existing rows 138-177 down to 139-178.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new price record.
$ws.Cells.Item(138,1).Value  = 4
$ws.Cells.Item(138,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(138,3).Value  = "Los Lagos"
$ws.Cells.Item(138,4).Value  = 44551
$ws.Cells.Item(138,5).Value  = 10
$ws.Cells.Item(138,6).Value  = "Fruta"
$ws.Cells.Item(138,7).Value  = 100108
$ws.Cells.Item(138,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(138,9).Value  = 100108005
$ws.Cells.Item(138,10).Value = "Piña"
$ws.Cells.Item(138,11).Value = "Caramelo"
$ws.Cells.Item(138,12).Value = "Tercera"
$ws.Cells.Item(138,13).Value = 300
$ws.Cells.Item(138,14).Value = 20000
$ws.Cells.Item(138,15).Value = 21000
$ws.Cells.Item(138,16).Value = 20500
$ws.Cells.Item(138,17).Value = "$/caja 16 unidades"
$ws.Cells.Item(138,18).Value = "Ecuador"
$ws.Cells.Item(138,19).Value = 1281
$ws.Cells.Item(138,20).Value = 16
